$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 14; $r++) {
    $dst = $ws.Range("R$r")
    $src = $ws.Range("Q$r")
    $dst.NumberFormat = $src.NumberFormat
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Bold = $src.Font.Bold
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.VerticalAlignment = $src.VerticalAlignment
}

$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 33.299999999999997
$ws.Range("R6").Value = 38.299999999999997
$ws.Range("R7").Value = 31.7
$ws.Range("R8").Value = 98.7
$ws.Range("R9").Value = 157.19999999999999
$ws.Range("R10").Value = 24.9
$ws.Range("R11").Value = 38.4
$ws.Range("R12").Value = 15.1
$ws.Range("R13").Value = 14.6
$ws.Range("R14").Value = 21.7

Write-Output "done"
